$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.096.43"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.651.72"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'533.83"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").Value = "'156.08"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'6.63"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "3.114.39"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "61.066.07"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "'22.12"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("D17").Value = "2.655.76"
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'355.74"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'10.69"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'61.72"
$ws.Range("E23").Value = "  +1.86%  "
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +6.89%  "
$ws.Range("E31").Value = "  +4.20%  "
$ws.Range("D32").Value = "'19.56"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'150.22"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'4.14"
$ws.Range("E34").Value = "  +3.23%  "
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "'0.925"
$ws.Range("E36").Value = "  +8.99%  "
$ws.Range("D37").Value = "'0.902"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").Value = "'310.61"
$ws.Range("E38").Value = "  +5.48%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'0.648"
$ws.Range("E41").Value = "  +3.56%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.99"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'5.01"
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").Value = "'19.26"
$ws.Range("E48").Value = "  +7.80%  "
$ws.Range("D49").Value = "'10.35"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "1.996.67"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "'1.85"
$ws.Range("E51").Value = "  +2.31%  "
